# Refresh the "data" sheet's time_taken column (F) with updated query timestamps
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value = "2021-10-05 14:21:37.671373"
$ws.Range("F3").Value = "2021-10-05 14:21:37.671380"
$ws.Range("F4").Value = "2021-10-05 14:21:37.671384"
$ws.Range("F5").Value = "2021-10-05 14:21:37.671386"
$ws.Range("F6").Value = "2021-10-05 14:21:37.671389"
$ws.Range("F7").Value = "2021-10-05 14:21:37.671392"
$ws.Range("F8").Value = "2021-10-05 14:21:37.671394"
$ws.Range("F9").Value = "2021-10-05 14:21:37.671397"
$ws.Range("F10").Value = "2021-10-05 14:21:37.671400"
$ws.Range("F11").Value = "2021-10-05 14:21:37.671402"
$ws.Range("F12").Value = "2021-10-05 14:21:37.671405"
$ws.Range("F13").Value = "2021-10-05 14:21:37.671408"
$ws.Range("F14").Value = "2021-10-05 14:21:37.671410"
$ws.Range("F15").Value = "2021-10-05 14:21:37.671413"
$ws.Range("F16").Value = "2021-10-05 14:21:37.671415"
$ws.Range("F17").Value = "2021-10-05 14:21:37.671418"
$ws.Range("F18").Value = "2021-10-05 14:21:37.671421"
$ws.Range("F19").Value = "2021-10-05 14:21:37.671423"
$ws.Range("F20").Value = "2021-10-05 14:21:37.671426"
$ws.Range("F21").Value = "2021-10-05 14:21:37.671428"
$ws.Range("F22").Value = "2021-10-05 14:21:37.671431"
$ws.Range("F23").Value = "2021-10-05 14:21:37.671433"
$ws.Range("F24").Value = "2021-10-05 14:21:37.671436"
$ws.Range("F25").Value = "2021-10-05 14:21:37.671439"
$ws.Range("F26").Value = "2021-10-05 14:21:37.671442"
$ws.Range("F27").Value = "2021-10-05 14:21:37.671444"
$ws.Range("F28").Value = "2021-10-05 14:21:37.671447"
$ws.Range("F29").Value = "2021-10-05 14:21:37.671450"
$ws.Range("F30").Value = "2021-10-05 14:21:37.671452"
$ws.Range("F31").Value = "2021-10-05 14:21:37.671455"
$ws.Range("F32").Value = "2021-10-05 14:21:37.671457"
$ws.Range("F33").Value = "2021-10-05 14:21:37.671460"
$ws.Range("F34").Value = "2021-10-05 14:21:37.671463"
$ws.Range("F35").Value = "2021-10-05 14:21:37.671466"
$ws.Range("F36").Value = "2021-10-05 14:21:37.671468"
$ws.Range("F37").Value = "2021-10-05 14:21:37.671471"
$ws.Range("F38").Value = "2021-10-05 14:21:37.671473"
$ws.Range("F39").Value = "2021-10-05 14:21:37.671476"
$ws.Range("F40").Value = "2021-10-05 14:21:37.671478"
$ws.Range("F41").Value = "2021-10-05 14:21:37.671481"
$ws.Range("F42").Value = "2021-10-05 14:21:37.671484"
$ws.Range("F43").Value = "2021-10-05 14:21:37.671487"
$ws.Range("F44").Value = "2021-10-05 14:21:37.671489"
$ws.Range("F45").Value = "2021-10-05 14:21:37.671492"
$ws.Range("F46").Value = "2021-10-05 14:21:37.671494"
$ws.Range("F47").Value = "2021-10-05 14:21:37.671497"

# Add a new "metadata" worksheet right after "data": duplicate "data" (so sheet-level
# settings like outline/page-setup match) then wipe its contents and rename it.
$ws.Copy($null, $ws)
$newSheet = $wb.Worksheets.Item("data (2)")
$newSheet.Cells.Clear()
$newSheet.Name = "metadata"

# Header row (B1:G1) - copy formatting (bold/border/center) from the data sheet's
# header row so the new sheet reuses the same style index, then overwrite the text.
$ws.Range("B1:F1").Copy($newSheet.Range("B1:F1"))
$ws.Range("B1").Copy($newSheet.Range("G1"))

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row 2 - copy A2's style (bold/border/center) from the data sheet, then set values
$ws.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").Value = "Mosaic skin disorders - deep sequencing"
$newSheet.Range("C2").Value = 564
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.5"
$newSheet.Range("E2").Value = "2021-01-06T15:23:02.528142Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:37.667773"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/564/?format=json"

# Restore "data" as the active sheet (tab selection unchanged in the source diff)
$ws.Activate()

Write-Output "metadata sheet added; timestamps refreshed"
